# Atualização de casos/óbitos até 28/06/22
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows (date serial, epidemiological_week col B, last_available_confirmed,
# last_available_deaths, new_confirmed, new_deaths) appended after existing row 65.
$data = @(
    @(44712, 0, 327458, 6348,   7, 0),
    @(44713, 0, 327474, 6348,  16, 0),
    @(44714, 0, 327490, 6348,  16, 0),
    @(44715, 0, 327517, 6348,  27, 0),
    @(44716, 0, 327528, 6348,  11, 0),
    @(44717, 0, 327538, 6348,  10, 0),
    @(44718, 0, 327552, 6348,  14, 0),
    @(44719, 0, 327574, 6348,  22, 0),
    @(44720, 0, 327597, 6349,  23, 1),
    @(44721, 0, 327638, 6349,  41, 0),
    @(44722, 0, 327687, 6350,  49, 1),
    @(44725, 0, 327737, 6351,  15, 0),
    @(44726, 0, 327777, 6352,  40, 1),
    @(44727, 0, 327910, 6353, 133, 1),
    @(44728, 0, 327993, 6353,  83, 0),
    @(44729, 0, 328042, 6353,  49, 0),
    @(44732, 0, 328119, 6353,  23, 0),
    @(44733, 0, 328260, 6354, 141, 1),
    @(44734, 0, 328561, 6354, 301, 0),
    @(44735, 0, 328744, 6355, 183, 1),
    @(44736, 0, 328955, 6355, 211, 0),
    @(44739, 0, 329405, 6356, 114, 0),
    @(44740, 0, 329686, 6357, 282, 1)
)

$startRow = 66
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 1).NumberFormat = "yyyy\-mm\-dd;@"
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
}

# Update the sheet view state to match the new scroll/selection position.
$win = $excel.ActiveWindow
$win.ScrollRow = 68
$win.ScrollColumn = 1
$ws.Range("E87").Select()
